# Updates the mounted-pipeline recompute: Step1_Data raw signal values for
# segment rows 6 and 11 (columns AQ:BW) changed upstream, which ripples into
# the Step2_Sj cumulative sums and the Step3_DataPts_* threshold-crossing
# lookups (columns D/F/G) for thresholds 0.5, 0.7, 0.8 and 0.9.
$wb = $excel.ActiveWorkbook

# --- Step1_Data: raw per-segment signal values (rows 6 and 11) ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("AQ6").Value = 0
$ws.Range("AR6").Value = 0.223246777998668
$ws.Range("AS6").Value = 0.1347025873881126
$ws.Range("AT6").Value = 0.2181735629871719
$ws.Range("AU6").Value = 0.01174227507609177
$ws.Range("AV6").Value = 0.009195967321217115
$ws.Range("AW6").Value = 0.0005034000952110847
$ws.Range("AX6").Value = 0.1265407682587021
$ws.Range("AY6").Value = 0.05617059619053
$ws.Range("AZ6").Value = 0.05308185065959252
$ws.Range("BA6").Value = 0.003202036977778243
$ws.Range("BB6").Value = 0.003099623613695048
$ws.Range("BC6").Value = 0.006534066280422819
$ws.Range("BD6").Value = 0.001063958216731477
$ws.Range("BE6").Value = 0.004403899190271915
$ws.Range("BF6").Value = 0.005908773701796164
$ws.Range("BG6").Value = 0.0004375472663692388
$ws.Range("BH6").Value = 0.01578577610639678
$ws.Range("BI6").Value = 0.006704223781297329
$ws.Range("BJ6").Value = 0.0003381909325252729
$ws.Range("BK6").Value = 0.02883040901600405
$ws.Range("BL6").Value = 0.03240346678072314
$ws.Range("BM6").Value = 0.002475035819328574
$ws.Range("BN6").Value = 0.01421128806915033
$ws.Range("BO6").Value = 0.0008569314440542653
$ws.Range("BP6").Value = 0.005954239813400096
$ws.Range("BQ6").Value = 0.003446600309258671
$ws.Range("BR6").Value = 0.0005414345190897467
$ws.Range("BS6").Value = 0.006211289686732841
$ws.Range("BT6").Value = 0.003088568536573119
$ws.Range("BU6").Value = 0.002111207253124406
$ws.Range("BV6").Value = 0.0006038142897018868
$ws.Range("BW6").Value = 0.0184298324202772
$ws.Range("AQ11").Value = 0
$ws.Range("AR11").Value = 0.2945010927386066
$ws.Range("AS11").Value = 0.09771066053687859
$ws.Range("AT11").Value = 0.2057252104147719
$ws.Range("AU11").Value = 0.01557321129169069
$ws.Range("AV11").Value = 0.00352068767261901
$ws.Range("AW11").Value = 0.002273184648482908
$ws.Range("AX11").Value = 0.07843114332949198
$ws.Range("AY11").Value = 0.02072393410740649
$ws.Range("AZ11").Value = 0.05774299577384605
$ws.Range("BA11").Value = 0.0108507835304889
$ws.Range("BB11").Value = 0.00002708553366157931
$ws.Range("BC11").Value = 0.01677939586542765
$ws.Range("BD11").Value = 0.01257204397400565
$ws.Range("BE11").Value = 0.03515048593735561
$ws.Range("BF11").Value = 0.008932756234424303
$ws.Range("BG11").Value = 0.0007025625441254645
$ws.Range("BH11").Value = 0.01718763237011445
$ws.Range("BI11").Value = 0.002896729489078782
$ws.Range("BJ11").Value = 0.002594886153465553
$ws.Range("BK11").Value = 0.0424554230967461
$ws.Range("BL11").Value = 0.02249858965191356
$ws.Range("BM11").Value = 0.0001056410992753355
$ws.Range("BN11").Value = 0.002250205809783782
$ws.Range("BO11").Value = 0.0003939916443419537
$ws.Range("BP11").Value = 0.0008161860692180302
$ws.Range("BQ11").Value = 0.002242357746053525
$ws.Range("BR11").Value = 0.0002339532016064927
$ws.Range("BS11").Value = 0.005944883694113478
$ws.Range("BT11").Value = 0.001166639490389394
$ws.Range("BU11").Value = 0.003642396199683159
$ws.Range("BV11").Value = 0.004015248078527668
$ws.Range("BW11").Value = 0.0303380020724051

# --- Step2_Sj: cumulative sum of Step1_Data (rows 6 and 11) ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("AQ6").Value = 0
$ws.Range("AR6").Value = 0.223246777998668
$ws.Range("AS6").Value = 0.3579493653867806
$ws.Range("AT6").Value = 0.5761229283739526
$ws.Range("AU6").Value = 0.5878652034500444
$ws.Range("AV6").Value = 0.5970611707712615
$ws.Range("AW6").Value = 0.5975645708664726
$ws.Range("AX6").Value = 0.7241053391251748
$ws.Range("AY6").Value = 0.7802759353157048
$ws.Range("AZ6").Value = 0.8333577859752973
$ws.Range("BA6").Value = 0.8365598229530755
$ws.Range("BB6").Value = 0.8396594465667706
$ws.Range("BC6").Value = 0.8461935128471935
$ws.Range("BD6").Value = 0.8472574710639249
$ws.Range("BE6").Value = 0.8516613702541969
$ws.Range("BF6").Value = 0.8575701439559931
$ws.Range("BG6").Value = 0.8580076912223623
$ws.Range("BH6").Value = 0.8737934673287591
$ws.Range("BI6").Value = 0.8804976911100564
$ws.Range("BJ6").Value = 0.8808358820425817
$ws.Range("BK6").Value = 0.9096662910585858
$ws.Range("BL6").Value = 0.9420697578393089
$ws.Range("BM6").Value = 0.9445447936586375
$ws.Range("BN6").Value = 0.9587560817277878
$ws.Range("BO6").Value = 0.9596130131718421
$ws.Range("BP6").Value = 0.9655672529852422
$ws.Range("BQ6").Value = 0.9690138532945008
$ws.Range("BR6").Value = 0.9695552878135906
$ws.Range("BS6").Value = 0.9757665775003234
$ws.Range("BT6").Value = 0.9788551460368965
$ws.Range("BU6").Value = 0.9809663532900209
$ws.Range("BV6").Value = 0.9815701675797228
$ws.Range("BW6").Value = 1
$ws.Range("AQ11").Value = 0
$ws.Range("AR11").Value = 0.2945010927386066
$ws.Range("AS11").Value = 0.3922117532754852
$ws.Range("AT11").Value = 0.5979369636902571
$ws.Range("AU11").Value = 0.6135101749819477
$ws.Range("AV11").Value = 0.6170308626545667
$ws.Range("AW11").Value = 0.6193040473030497
$ws.Range("AX11").Value = 0.6977351906325416
$ws.Range("AY11").Value = 0.7184591247399481
$ws.Range("AZ11").Value = 0.7762021205137941
$ws.Range("BA11").Value = 0.787052904044283
$ws.Range("BB11").Value = 0.7870799895779446
$ws.Range("BC11").Value = 0.8038593854433722
$ws.Range("BD11").Value = 0.8164314294173779
$ws.Range("BE11").Value = 0.8515819153547335
$ws.Range("BF11").Value = 0.8605146715891578
$ws.Range("BG11").Value = 0.8612172341332833
$ws.Range("BH11").Value = 0.8784048665033978
$ws.Range("BI11").Value = 0.8813015959924766
$ws.Range("BJ11").Value = 0.8838964821459422
$ws.Range("BK11").Value = 0.9263519052426883
$ws.Range("BL11").Value = 0.9488504948946018
$ws.Range("BM11").Value = 0.9489561359938772
$ws.Range("BN11").Value = 0.9512063418036609
$ws.Range("BO11").Value = 0.9516003334480029
$ws.Range("BP11").Value = 0.9524165195172209
$ws.Range("BQ11").Value = 0.9546588772632745
$ws.Range("BR11").Value = 0.9548928304648809
$ws.Range("BS11").Value = 0.9608377141589944
$ws.Range("BT11").Value = 0.9620043536493839
$ws.Range("BU11").Value = 0.9656467498490671
$ws.Range("BV11").Value = 0.9696619979275948
$ws.Range("BW11").Value = 0.9999999999999999

# --- Step3_DataPts_0.5: threshold-crossing point for rows 6 and 11 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F6").Value = 0.5761229283739526
$ws.Range("F11").Value = 0.5979369636902571

# --- Step3_DataPts_0.7: threshold-crossing point for rows 6 and 11 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F6").Value = 0.7241053391251748
$ws.Range("D11").Value = 50
$ws.Range("F11").Value = 0.7184591247399481
$ws.Range("G11").Value = 9

# --- Step3_DataPts_0.8: threshold-crossing point for rows 6 and 11 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F6").Value = 0.8333577859752973
$ws.Range("D11").Value = 54
$ws.Range("F11").Value = 0.8038593854433722
$ws.Range("G11").Value = 13

# --- Step3_DataPts_0.9: threshold-crossing point for rows 6 and 11 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F6").Value = 0.9096662910585858
$ws.Range("D11").Value = 62
$ws.Range("F11").Value = 0.9263519052426883
$ws.Range("G11").Value = 21
